# Rename the embedded logo pictures' docPr/name identifiers.
#
#   Footer (first page)  - id=3 Pearson logo : image1.png -> image2.png
#   Footer (default/odd) - id=2 Pearson logo : image1.png -> image2.png
#   Header (first page)  - id=1 BTec logo    : image2.jpg -> image1.jpg
#
# wdHeaderFooterIndex: 1 = primary/default, 2 = first page, 3 = even pages.

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# --- Footers -------------------------------------------------------------
# Footer index 2 -> word/footer1.xml (first page) -- Pearson logo docPr id="3"
$footerFirst = $section.Footers.Item(2)
if ($footerFirst.Exists) {
    $logo = $footerFirst.Range.InlineShapes.Item(1)
    $logo.Name = "image2.png"
}

# Footer index 1 -> word/footer2.xml (default/odd pages) -- Pearson logo docPr id="2"
$footerDefault = $section.Footers.Item(1)
if ($footerDefault.Exists) {
    $logo = $footerDefault.Range.InlineShapes.Item(1)
    $logo.Name = "image2.png"
}

# --- Headers ---------------------------------------------------------------
# Header index 2 -> word/header1.xml (first page) -- BTec logo docPr id="1"
$headerFirst = $section.Headers.Item(2)
if ($headerFirst.Exists) {
    $logo = $headerFirst.Range.InlineShapes.Item(1)
    $logo.Name = "image1.jpg"
}
